$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.ListFormat.ListType -ne 0) {
        $p.Range.ListFormat.RemoveNumbers()
    }
}
